# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1, style index 1:
# bold, bordered, centered) onto the new header cell H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Header text for the new column.
$ws.Range("H1").Value = "Save"

# Data value for the new column (row 2).
$ws.Range("H2").Value = 1
